# Chauffeurs.xlsx - remove duplicate/stale driver rows and fix a few
# lastname/firstname entries that were entered with the wrong surname.
#
# "Import works (kind of)" - cleans up a handful of rows that were
# imported with merged/incorrect surnames, and drops rows that turned
# out to be duplicates no longer needed.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the 5 rows that are no longer needed, bottom-to-top so the
# row numbers of the ones still to be removed don't shift underneath us.
$ws.Rows.Item(56).Delete()   # Ramella, Elodie
$ws.Rows.Item(45).Delete()   # Marquis, Marie Jo
$ws.Rows.Item(44).Delete()   # Mariotti, Monica
$ws.Rows.Item(35).Delete()   # Janin Cancian, Léonore
$ws.Rows.Item(30).Delete()   # Glaus, Vania

# Fix a handful of lastname / firstname values that were imported with
# the wrong surname.
$ws.Range("A3").Value = "Angiolili"
$ws.Range("A20").Value = "Féry-Hammer"
$ws.Range("A22").Value = "Fleischman"
$ws.Range("A49").Value = "Pinilla-Marin"
$ws.Range("B49").Value = "Andres"

# Re-apply the AutoFilter over the now-smaller range (A1:D62) instead of
# the old A1:D67, and keep the hidden _FilterDatabase name in sync.
$ws.AutoFilterMode = $false
$ws.Range("A1:D62").AutoFilter()

$n = $wb.Names.Item(0)
$n.RefersTo = "=Chauffeurs!`$A`$1:`$D`$62"

# Move the active selection (no more scrolled-down view / old selection).
$ws.Range("K75").Select()
